# Trade #13 closed at 2026-02-18 00:10:50 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1499.73
$wsSummary.Range("B4").Value = 0.83
$wsSummary.Range("B5").Value = 0.4
$wsSummary.Range("B6").Value = 41
$wsSummary.Range("B7").Value = 23
$wsSummary.Range("B9").Value = 56.1

# --- Strategy Status sheet (MarketMaking row, row 6) ---
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C6").Value = 99.73
$wsStatus.Range("D6").Value = 12
$wsStatus.Range("E6").Value = -0.08
$wsStatus.Range("F6").Value = -0.27
$wsStatus.Range("G6").Value = 58.33

# --- All Trades sheet (trade row 43, trade #42 closing) ---
$wsTrades = $wb.Worksheets.Item("All Trades")
$wsTrades.Range("G43").Value = 0.44
$wsTrades.Range("H43").Value = "CLOSED"
$wsTrades.Range("I43").Value = 4.7619
$wsTrades.Range("J43").Value = 0.02
$wsTrades.Range("K43").Value = 99.73
$wsTrades.Range("L43").Value = "early_exit"
$wsTrades.Range("M43").Value = 0.13

# --- MarketMaking sheet (trade row 14, trade #42 closing) ---
$wsMM = $wb.Worksheets.Item("MarketMaking")
$wsMM.Range("G14").Value = 0.44
$wsMM.Range("H14").Value = "CLOSED"
$wsMM.Range("I14").Value = 4.7619
$wsMM.Range("J14").Value = 0.02
$wsMM.Range("K14").Value = 99.73
$wsMM.Range("P14").Value = "early_exit"
$wsMM.Range("Q14").Value = 0.13
